$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.NumberFormat = "General"
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '24.749.17'
$ws.Range("E2").Value = '  +0.48%  '

Set-TextValue $ws.Range("D3") '1.703.00'
$ws.Range("E3").Value = '  +0.43%  '

Set-TextValue $ws.Range("D4") '1.004'
$ws.Range("E4").Value = '  +0.25%  '

Set-TextValue $ws.Range("D5") '317.37'
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("E6").Value = '  +0.24%  '

Set-TextValue $ws.Range("D7") '0.3950'
$ws.Range("E7").Value = '  -0.19%  '

Set-TextValue $ws.Range("D8") '0.4058'
$ws.Range("E8").Value = '  +0.86%  '

Set-TextValue $ws.Range("D9") '1.537'
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("E10").Value = '  +0.19%  '

Set-TextValue $ws.Range("D11") '53.63'
$ws.Range("E11").Value = '  -0.32%  '

Set-TextValue $ws.Range("D12") '0.08900'
$ws.Range("E12").Value = '  +1.17%  '

Set-TextValue $ws.Range("D13") '7.465'
$ws.Range("E13").Value = '  +2.83%  '

Set-TextValue $ws.Range("D14") '23.78'
$ws.Range("E14").Value = '  +2.02%  '

Set-TextValue $ws.Range("D15") '8.183'
$ws.Range("E15").Value = '  +7.25%  '

Set-TextValue $ws.Range("D16") '0.00001327'
$ws.Range("E16").Value = '  +0.49%  '

Set-TextValue $ws.Range("D17") '1.708.57'
$ws.Range("E17").Value = '  +0.38%  '

Set-TextValue $ws.Range("D18") '99.91'
$ws.Range("E18").Value = '  -1.28%  '

Set-TextValue $ws.Range("D19") '0.07068'
$ws.Range("E19").Value = '  +0.71%  '

Set-TextValue $ws.Range("D20") '19.77'
$ws.Range("E20").Value = '  +0.23%  '

Set-TextValue $ws.Range("D21") '7.091'
$ws.Range("E21").Value = '  +2.83%  '

$ws.Range("E22").Value = '  +0.14%  '

Set-TextValue $ws.Range("D23") '14.65'
$ws.Range("E23").Value = '  +3.88%  '

Set-TextValue $ws.Range("D24") '24.732.56'
$ws.Range("E24").Value = '  +0.47%  '

Set-TextValue $ws.Range("D25") '3.226'
$ws.Range("E25").Value = '  +5.15%  '

Set-TextValue $ws.Range("D26") '2.374'
$ws.Range("E26").Value = '  +1.59%  '

Set-TextValue $ws.Range("D27") '22.82'
$ws.Range("E27").Value = '  +2.04%  '

Set-TextValue $ws.Range("D28") '162.93'
$ws.Range("E28").Value = '  +1.89%  '

Set-TextValue $ws.Range("D29") '8.782'
$ws.Range("E29").Value = '  +18.09%  '

Set-TextValue $ws.Range("D30") '136.26'
$ws.Range("E30").Value = '  +1.43%  '

$ws.Range("E31").Value = '  -1.40%  '

Set-TextValue $ws.Range("D32") '7.731'
$ws.Range("E32").Value = '  +3.77%  '

Set-TextValue $ws.Range("D33") '0.08959'
$ws.Range("E33").Value = '  +4.97%  '

Set-TextValue $ws.Range("D34") '1.080'
$ws.Range("E34").Value = '  -3.08%  '

Set-TextValue $ws.Range("D35") '1.994'
$ws.Range("E35").Value = '  +1.07%  '

Set-TextValue $ws.Range("D36") '11.13'
$ws.Range("E36").Value = '  -2.92%  '

Set-TextValue $ws.Range("D37") '0.2765'
$ws.Range("E37").Value = '  +0.65%  '

Set-TextValue $ws.Range("D38") '14.59'
$ws.Range("E38").Value = '  -0.29%  '

Set-TextValue $ws.Range("D39") '0.02791'
$ws.Range("E39").Value = '  +0.08%  '

Set-TextValue $ws.Range("D40") '0.09191'
$ws.Range("E40").Value = '  +1.52%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D41") '0.7740'
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D42") '1.461'
$ws.Range("E42").Value = '  -0.38%  '

Set-TextValue $ws.Range("D43") '16.10'
$ws.Range("E43").Value = '  +3.96%  '

Set-TextValue $ws.Range("D44") '0.7222'
$ws.Range("E44").Value = '  -0.18%  '

Set-TextValue $ws.Range("D45") '2.582'
$ws.Range("E45").Value = '  +2.44%  '

Set-TextValue $ws.Range("D46") '4.226'
$ws.Range("E46").Value = '  -0.25%  '

Set-TextValue $ws.Range("D47") '1.367'
$ws.Range("E47").Value = '  +0.21%  '

Set-TextValue $ws.Range("D48") '1.004'
$ws.Range("E48").Value = '  +0.27%  '

Set-TextValue $ws.Range("D49") '140.84'
$ws.Range("E49").Value = '  -0.14%  '

Set-TextValue $ws.Range("D50") '91.18'
$ws.Range("E50").Value = '  +2.98%  '

$ws.Range("E51").Value = '  -0.60%  '
